{"js": "// Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n// (percentages, dollar amounts, large numbers) across the resume body.\n//\n// Strategy: for each target paragraph, locate each metric substring with\n// `paragraph.search(...)` (scoped to that paragraph only, so repeats of the\n// same literal in other paragraphs are unaffected) and apply\n// bold + color (#2C3E50) to the matched run. Word's search naturally splits\n// the original single run into \"before\" / \"match\" / \"after\" runs, which is\n// exactly the run structure the diff shows.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\nconst SEARCH_OPTIONS = { matchCase: true, matchWholeWord: false };\n\n// Map of a distinctive paragraph-prefix (enough to uniquely identify the\n// paragraph without depending on a fixed index) -> ordered list of metric\n// substrings to bold+color within that paragraph.\nconst targets = [\n  {\n    prefix: \"\u2022 Discovered systematic race coding errors\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    prefix: \"\u2022 Utilized advanced sampling methods\",\n    metrics: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    prefix: \"\u2022 Trigonometric algorithm for boundary estimation\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    prefix: \"\u2022 Built real-time FEC analysis systems\",\n    metrics: [\"$2\"],\n  },\n  {\n    prefix: \"\u2022 Modernized legacy ETL processes\",\n    metrics: [\"57%\"],\n  },\n  {\n    prefix: \"\u2022 Platform impact: Built redistricting system serving\",\n    metrics: [\"12,847\"],\n  },\n  {\n    prefix: \"\u2022 Revenue generation: Delivered\",\n    metrics: [\"$4.9M\"],\n  },\n  {\n    prefix: \"\u2022 23% conversion rate improvement\",\n    metrics: [\"23%\"],\n  },\n];\n\nfor (const target of targets) {\n  const para = paragraphs.items.find((p) => p.text.indexOf(target.prefix) === 0);\n  if (!para) {\n    continue;\n  }\n  for (const metric of target.metrics) {\n    const found = para.search(metric, SEARCH_OPTIONS);\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length === 0) {\n      continue;\n    }\n    // Only the first occurrence is the intended metric (each literal is\n    // unique within its own paragraph for this document).\n    const matchRange = found.items[0];\n    matchRange.font.set({ bold: true, color: HIGHLIGHT_COLOR });\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n# (percentages, dollar amounts, large numbers) across the resume body.\n#\n# Strategy: for each target paragraph (identified by a distinctive literal\n# prefix), walk its metrics in left-to-right order. For each metric, run\n# Range.Find.Execute scoped to a shrinking sub-range (from the end of the\n# previous match through the end of the paragraph) so repeated literals\n# elsewhere in the document/paragraph are never touched, then bold + color\n# (#2C3E50) the matched run. Word's Find naturally splits the original\n# single run into \"before\" / \"match\" / \"after\" runs, matching the diff.\n\n$d = $word.ActiveDocument\n\n# RGB(0x2C, 0x3E, 0x50) packed as VBA/COM expects (R + G*256 + B*65536).\n$highlightColor = 0x2C + (0x3E * 256) + (0x50 * 65536)\n\n# Ordered list of (paragraph-prefix, metrics-in-order) pairs, mirroring the\n# commit's regex-driven highlighter applied to achievements/responsibilities.\n$targets = @(\n    @(\"\u2022 Discovered systematic race coding errors\", @(\"23%\", \"64%\")),\n    @(\"\u2022 Utilized advanced sampling methods\", @([char]0x00B1 + \"4.2%\", [char]0x00B1 + \"2.1%\", \"71%\", \"87%\")),\n    @(\"\u2022 Trigonometric algorithm for boundary estimation\", @(\"73.5%\", \"$4.7M\")),\n    @(\"\u2022 Built real-time FEC analysis systems\", @(\"$2\")),\n    @(\"\u2022 Modernized legacy ETL processes\", @(\"57%\")),\n    @(\"\u2022 Platform impact: Built redistricting system serving\", @(\"12,847\")),\n    @(\"\u2022 Revenue generation: Delivered\", @(\"$4.9M\")),\n    @(\"\u2022 23% conversion rate improvement\", @(\"23%\"))\n)\n\nforeach ($entry in $targets) {\n    $prefix = $entry[0]\n    $metrics = $entry[1]\n\n    $target = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.StartsWith($prefix)) {\n            $target = $p\n        }\n    }\n    if ($target -eq $null) {\n        continue\n    }\n\n    $paraRange = $target.Range\n    $cursor = $paraRange.Start\n    $paraEnd = $paraRange.End\n\n    foreach ($metric in $metrics) {\n        $searchRange = $d.Range($cursor, $paraEnd)\n        $found = $searchRange.Find.Execute($metric, $true, $false)\n        if (-not $found) {\n            continue\n        }\n        $searchRange.Font.Bold = 1\n        $searchRange.Font.Color = $highlightColor\n        $cursor = $searchRange.End\n    }\n}\n"}
